$wb = $excel.ActiveWorkbook

$renters = $wb.Worksheets.Item("renters")
$owners  = $wb.Worksheets.Item("owners")

# --- renters sheet: zero-out savings column L for rows 2-5 ---
$renters.Range("L2").Value = 0
$renters.Range("L3").Value = 0
$renters.Range("L4").Value = 0
$renters.Range("L5").Value = 0

# --- owners sheet: zero-out savings column B for rows 2-5 ---
$owners.Range("B2").Value = 0
$owners.Range("B3").Value = 0
$owners.Range("B4").Value = 0
$owners.Range("B5").Value = 0

# --- owners sheet: reset damage state column L to "Slight" ---
$owners.Range("L2").Value = "Slight"
$owners.Range("L3").Value = "Slight"
$owners.Range("L4").Value = "Slight"
$owners.Range("L5").Value = "Slight"
$owners.Range("L7").Value = "Slight"
$owners.Range("L8").Value = "Slight"
$owners.Range("L9").Value = "Slight"

# --- owners sheet: update owner credit column P row 2 ---
$owners.Range("P2").Value = 700

# --- selections: move away from renters, make owners the active/selected sheet ---
$renters.Range("Q20").Select()
$owners.Range("L2:L9").Select()
$owners.Activate()
